$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 316354
$ws.Range("D2").Value = 403223042
$ws.Range("C8").Value = 851
$ws.Range("D8").Value = 1251908
$ws.Range("C10").Value = 116093
$ws.Range("D10").Value = 170114598
$ws.Range("C12").Value = 58616
$ws.Range("D12").Value = 84594215
$ws.Range("C16").Value = 3972
$ws.Range("D16").Value = 5637019
$ws.Range("C20").Value = 6493
$ws.Range("D20").Value = 9054961
$ws.Range("C22").Value = 76465
$ws.Range("D22").Value = 95409357
$ws.Range("C28").Value = 32219
$ws.Range("D28").Value = 47169923
$ws.Range("C30").Value = 11360
$ws.Range("D30").Value = 16339647
$ws.Range("C33").Value = 1555
$ws.Range("D33").Value = 2183807
$ws.Range("C35").Value = 1782
$ws.Range("D35").Value = 2514833
$ws.Range("C36").Value = 96056
$ws.Range("D36").Value = 120972935
$ws.Range("C42").Value = 899
$ws.Range("D42").Value = 1323185
$ws.Range("C44").Value = 44061
$ws.Range("D44").Value = 64575526
$ws.Range("C46").Value = 9036
$ws.Range("D46").Value = 12969198
$ws.Range("C48").Value = 1395
$ws.Range("D48").Value = 1937311
$ws.Range("C51").Value = 2257
$ws.Range("D51").Value = 3149001
$ws.Range("C52").Value = 68194
$ws.Range("D52").Value = 85580290
$ws.Range("C58").Value = 27916
$ws.Range("D58").Value = 40942029
$ws.Range("C61").Value = 10958
$ws.Range("D61").Value = 15843956
$ws.Range("C67").Value = 1441
$ws.Range("D67").Value = 2017921
$ws.Range("C69").Value = 20252
$ws.Range("D69").Value = 26526763
$ws.Range("C73").Value = 7516
$ws.Range("D73").Value = 11003530
$ws.Range("C75").Value = 5053
$ws.Range("D75").Value = 7337606
$ws.Range("C77").Value = 269
$ws.Range("D77").Value = 377173
$ws.Range("C78").Value = 138960
$ws.Range("D78").Value = 173322913
$ws.Range("C84").Value = 63002
$ws.Range("D84").Value = 92342999
$ws.Range("C87").Value = 29370
$ws.Range("D87").Value = 42485087
$ws.Range("C89").Value = 2717
$ws.Range("D89").Value = 3912860
$ws.Range("C90").Value = 2764
$ws.Range("D90").Value = 3906280
$ws.Range("C91").Value = 32126
$ws.Range("D91").Value = 43525558
$ws.Range("C95").Value = 7802
$ws.Range("D95").Value = 11472470
$ws.Range("C97").Value = 7098
$ws.Range("D97").Value = 10290127
$ws.Range("C101").Value = 8823
$ws.Range("D101").Value = 12245320
$ws.Range("C103").Value = 2218
$ws.Range("D103").Value = 3267382
$ws.Range("C105").Value = 2988
$ws.Range("D105").Value = 4364741
$ws.Range("C108").Value = 175
$ws.Range("D108").Value = 248586
$ws.Range("C109").Value = 139425
$ws.Range("D109").Value = 172424855
$ws.Range("C111").Value = 72
$ws.Range("D111").Value = 102644
$ws.Range("C115").Value = 52232
$ws.Range("D115").Value = 76569588
$ws.Range("C117").Value = 26624
$ws.Range("D117").Value = 38572541
$ws.Range("C118").Value = 1301
$ws.Range("D118").Value = 1780991
$ws.Range("C121").Value = 2203
$ws.Range("D121").Value = 3093455
$ws.Range("C123").Value = 494235
$ws.Range("D123").Value = 651778271
$ws.Range("C130").Value = 204955
$ws.Range("D130").Value = 301289885
$ws.Range("C133").Value = 177110
$ws.Range("D133").Value = 257427534
$ws.Range("C136").Value = 2820
$ws.Range("D136").Value = 3963285
$ws.Range("C138").Value = 6177
$ws.Range("D138").Value = 8727534
$ws.Range("C141").Value = 43856
$ws.Range("D141").Value = 58552963
$ws.Range("C147").Value = 13905
$ws.Range("D147").Value = 20395249
$ws.Range("C148").Value = 3702
$ws.Range("D148").Value = 5338019
$ws.Range("C151").Value = 395
$ws.Range("D151").Value = 568431
$ws.Range("C154").Value = 17264
$ws.Range("D154").Value = 22809991
$ws.Range("C158").Value = 7059
$ws.Range("D158").Value = 10266288
$ws.Range("C160").Value = 4913
$ws.Range("D160").Value = 7071877
$ws.Range("C165").Value = 15462
$ws.Range("D165").Value = 22435556
$ws.Range("C166").Value = 1744
$ws.Range("D166").Value = 2593730
$ws.Range("C171").Value = 86613
$ws.Range("D171").Value = 108352476
$ws.Range("C178").Value = 33549
$ws.Range("D178").Value = 49200257
$ws.Range("C180").Value = 12840
$ws.Range("D180").Value = 18550068
$ws.Range("C182").Value = 1239
$ws.Range("D182").Value = 1734396
$ws.Range("C184").Value = 1611
$ws.Range("D184").Value = 2263302
$ws.Range("C186").Value = 235414
$ws.Range("D186").Value = 292679928
$ws.Range("C192").Value = 866
$ws.Range("D192").Value = 1273997
$ws.Range("C194").Value = 85906
$ws.Range("D194").Value = 125928393
$ws.Range("C197").Value = 32653
$ws.Range("D197").Value = 46994962
$ws.Range("C200").Value = 5053
$ws.Range("D200").Value = 7201413
$ws.Range("C203").Value = 4756
$ws.Range("D203").Value = 6582573
$ws.Range("C206").Value = 260327
$ws.Range("D206").Value = 322224123
$ws.Range("C215").Value = 94314
$ws.Range("D215").Value = 137978886
$ws.Range("C218").Value = 50799
$ws.Range("D218").Value = 73415213
$ws.Range("C221").Value = 4637
$ws.Range("D221").Value = 6509957
$ws.Range("C224").Value = 5602
$ws.Range("D224").Value = 7746365
$ws.Range("C227").Value = 104843
$ws.Range("D227").Value = 131197664
$ws.Range("C234").Value = 49075
$ws.Range("D234").Value = 71896047
$ws.Range("C236").Value = 12221
$ws.Range("D236").Value = 17570469
$ws.Range("C240").Value = 2440
$ws.Range("D240").Value = 3409315
$ws.Range("C241").Value = 253871
$ws.Range("D241").Value = 320583029
$ws.Range("C242").Value = 171
$ws.Range("D242").Value = 211433
$ws.Range("C249").Value = 94875
$ws.Range("D249").Value = 139021967
$ws.Range("C252").Value = 64078
$ws.Range("D252").Value = 92856897
$ws.Range("C254").Value = 2389
$ws.Range("D254").Value = 3370428
$ws.Range("C257").Value = 4502
$ws.Range("D257").Value = 6320164
